$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, derived from the authoritative diff
$changes = @{
    2 = @{ D='63.919.28'; E='  +3.08%  ' }
    3 = @{ D='2.539.40'; E='  +5.54%  ' }
    4 = @{ D='1.00'; E='  +0.03%  ' }
    5 = @{ D='573.86'; E='  +2.23%  ' }
    6 = @{ D='147.79'; E='  +6.45%  ' }
    7 = @{ D='1.00'; E='  +0.06%  ' }
    8 = @{ D='0.590'; E='  +0.62%  ' }
    9 = @{ D='2.537.66'; E='  +5.54%  ' }
    10 = @{ D='0.107'; E='  +2.28%  ' }
    11 = @{ D='5.78'; E='  +1.13%  ' }
    12 = @{ E='  +1.73%  ' }
    13 = @{ D='0.360'; E='  +3.18%  ' }
    14 = @{ D='27.87'; E='  +8.70%  ' }
    15 = @{ D='2.993.47'; E='  +5.66%  ' }
    16 = @{ D='63.677.06'; E='  +2.82%  ' }
    17 = @{ D='0.0000145'; E='  +4.17%  ' }
    18 = @{ D='2.538.27'; E='  +5.42%  ' }
    19 = @{ D='11.52'; E='  +4.68%  ' }
    20 = @{ D='343.93'; E='  +0.64%  ' }
    21 = @{ D='4.36'; E='  +3.33%  ' }
    22 = @{ D='6.89'; E='  +0.05%  ' }
    23 = @{ E='  +0.57%  ' }
    24 = @{ D='66.12'; E='  +1.93%  ' }
    25 = @{ E='  -0.59%  ' }
    26 = @{ E='  +5.36%  ' }
    27 = @{ D='0.999'; E='  -0.07%  ' }
    28 = @{ D='8.26'; E='  -0.79%  ' }
    29 = @{ D='1.42'; E='  +3.11%  ' }
    30 = @{ D='0.0₃0827'; E='  +7.44%  ' }
    31 = @{ E='  +4.28%  ' }
    32 = @{ D='6.85'; E='  +7.73%  ' }
    33 = @{ D='177.06'; E='  +3.49%  ' }
    34 = @{ B='ImmutableX'; C='https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D='1.55'; E='  +11.25%  ' }
    35 = @{ B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='416.17'; E='  +14.42%  ' }
    36 = @{ D='0.405'; E='  +3.08%  ' }
    37 = @{ D='19.13'; E='  +3.57%  ' }
    38 = @{ D='4.45'; E='  -2.52%  ' }
    39 = @{ E='  -0.01%  ' }
    40 = @{ E='  +6.15%  ' }
    41 = @{ D='0.999'; E='  -0.04%  ' }
    42 = @{ D='40.78'; E='  +4.65%  ' }
    43 = @{ D='152.96' }
    44 = @{ D='3.80'; E='  +3.78%  ' }
    45 = @{ D='20.95'; E='  +2.65%  ' }
    46 = @{ E='  +4.30%  ' }
    47 = @{ D='0.0531'; E='  +2.77%  ' }
    48 = @{ E='  +0.69%  ' }
    49 = @{ D='19.04'; E='  +6.96%  ' }
    50 = @{ E='  +4.74%  ' }
    51 = @{ B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='1.84'; E='  +9.82%  ' }
}

foreach ($row in $changes.Keys) {
    $rowData = $changes[$row]
    foreach ($col in $rowData.Keys) {
        $addr = "$col$row"
        if ($col -eq 'D') {
            # Price column: force text so numeric-looking strings
            # (e.g. '1.00', '0.999') keep their exact printed form
            # instead of being normalized/rounded as numbers.
            $ws.Range($addr).Value = "'" + $rowData[$col]
        } else {
            $ws.Range($addr).Value = $rowData[$col]
        }
    }
}
